$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing header row (row 5), pushing
# everything down by 2 rows.
$ws.Rows.Item(5).Resize(2).Insert()

# Row 5: the new URL text, styled like the existing hyperlink-style cell (A4).
$ws.Range("A5").Value = "https://www.inkstain.net/2023/02/deadpool-diaries-ignore-this-post-about-the-latest-colorado-river-runoff-forecast/"
$ws.Range("A5").Style = "Hyperlink"

# Row 6: empty cell but carries the same Hyperlink style.
$ws.Range("A6").Style = "Hyperlink"

# Update the active selection to match the target state.
$ws.Range("B2:B3").Select()
